# Trade #21 closed at 2026-02-17 23:57:07 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.73
$summary.Range("B4").Value = 0.73
$summary.Range("B5").Value = 0.7
$summary.Range("B6").Value = 21
$summary.Range("B8").Value = 8
$summary.Range("B9").Value = 52.38

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.73
$status.Range("D6").Value = 21
$status.Range("E6").Value = 0.73
$status.Range("F6").Value = 0.73
$status.Range("G6").Value = 52.38

# ---------------------------------------------------------------------------
# Append new trade row (#21) to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
$newTradeRow = @(21, "2026-02-17", "23:57:00", "MarketMaking", "DOWN", 0.43, 0.4, "CLOSED", -6.9767, -0.03, 100.73, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.15)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 22
    $ws.Cells.Item($row, 1).Value = $newTradeRow[0]
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newTradeRow[1]
    $dateCell.Style = "Normal"
    $ws.Cells.Item($row, 3).Value = $newTradeRow[2]
    $ws.Cells.Item($row, 4).Value = $newTradeRow[3]
    $ws.Cells.Item($row, 5).Value = $newTradeRow[4]
    $ws.Cells.Item($row, 6).Value = $newTradeRow[5]
    $ws.Cells.Item($row, 7).Value = $newTradeRow[6]
    $ws.Cells.Item($row, 8).Value = $newTradeRow[7]
    $ws.Cells.Item($row, 9).Value = $newTradeRow[8]
    $ws.Cells.Item($row, 10).Value = $newTradeRow[9]
    $ws.Cells.Item($row, 11).Value = $newTradeRow[10]
    $ws.Cells.Item($row, 12).Value = $newTradeRow[11]
    $ws.Cells.Item($row, 13).Value = $newTradeRow[12]
    $ws.Cells.Item($row, 14).Value = $newTradeRow[13]
    $ws.Cells.Item($row, 15).Value = $newTradeRow[14]
    $ws.Cells.Item($row, 16).Value = $newTradeRow[15]
    $ws.Cells.Item($row, 17).Value = $newTradeRow[16]
}
